# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# to reflect the latest scrape snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new value for column F on sheet "展览"
$updatesExhibition = @{
    2  = 14892
    3  = 18590
    5  = 114
    14 = 112
    15 = 201
    17 = 1419
    18 = 151
    20 = 85
    22 = 7710
    28 = 5967
    29 = 104
    30 = 65
    33 = 258
    34 = 5326
}

# Row -> new value for column F on sheet "全部类型"
# (same events, but a couple of extra rows shift the numbering by +1
# starting partway down the sheet).
$updatesAllTypes = @{
    2  = 14892
    3  = 18590
    5  = 114
    14 = 112
    15 = 201
    17 = 1419
    18 = 151
    21 = 85
    23 = 7710
    31 = 5967
    32 = 104
    33 = 65
    36 = 258
    37 = 5326
}

$ws = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $ws.Cells.Item($row, 6).Value = $updatesExhibition[$row]
}

$ws = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAllTypes.Keys) {
    $ws.Cells.Item($row, 6).Value = $updatesAllTypes[$row]
}
